$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("model" row): correct typos / update the Description wording (added
# commas between the supported-cases bullets) and the Type column text.
$ws.Range("C2").Value = "Model path or model instance. The following cases are supported:" + [char]10 + "Using backend=""onnx"" and a onnx model path," + [char]10 + "Using backend=""tvm"" and a Keras model," + [char]10 + "Using backend=""tvm"" and a PyTorch model," + [char]10 + "Using backend=""tf"" and a fronzen TensorFlow PB graph"
# Leading apostrophe preserves the original quote-prefix (text) cell style.
$ws.Range("B2").Value = "'string or model class of the corresponding framework"

# Update the active cell selection
$ws.Range("B2").Select() | Out-Null
